$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2418952618453865
$ws.Range("C2").Value = 0.4663341645885287
$ws.Range("J2").Value = 0.009975062344139651
$ws.Range("P2").Value = 0.1895261845386534
$ws.Range("S2").Value = 0.09226932668329177
# Row 3
$ws.Range("B3").Value = 0.01063829787234043
$ws.Range("J3").Value = 0.0425531914893617
$ws.Range("P3").Value = 0.7074468085106383
$ws.Range("S3").Value = 0.2393617021276596
# Row 4
$ws.Range("J4").Value = 0.03529411764705882
$ws.Range("P4").Value = 0.6588235294117647
$ws.Range("S4").Value = 0.3058823529411765
# Row 6
$ws.Range("B6").Value = 0.06666666666666667
$ws.Range("D6").Value = 0.008333333333333333
$ws.Range("F6").Value = 0.05
$ws.Range("J6").Value = 0.2958333333333333
$ws.Range("O6").Value = 0.0375
$ws.Range("Q6").Value = 0.1625
$ws.Range("R6").Value = 0.0375
$ws.Range("S6").Value = 0.3416666666666667
# Row 7
$ws.Range("B7").Value = 0.1008064516129032
$ws.Range("D7").Value = 0.04435483870967742
$ws.Range("F7").Value = 0.03225806451612903
$ws.Range("J7").Value = 0.1411290322580645
$ws.Range("O7").Value = 0.02419354838709677
$ws.Range("Q7").Value = 0.1774193548387097
$ws.Range("R7").Value = 0.0846774193548387
$ws.Range("S7").Value = 0.3951612903225806
# Row 8
$ws.Range("B8").Value = 0.09172259507829977
$ws.Range("D8").Value = 0.029082774049217
$ws.Range("F8").Value = 0.06935123042505593
$ws.Range("J8").Value = 0.1319910514541387
$ws.Range("O8").Value = 0.02237136465324385
$ws.Range("Q8").Value = 0.174496644295302
$ws.Range("R8").Value = 0.04697986577181208
$ws.Range("S8").Value = 0.4340044742729307
# Row 9
$ws.Range("B9").Value = 0.1137254901960784
$ws.Range("D9").Value = 0.0196078431372549
$ws.Range("F9").Value = 0.04313725490196078
$ws.Range("J9").Value = 0.1137254901960784
$ws.Range("O9").Value = 0.02745098039215686
$ws.Range("Q9").Value = 0.2156862745098039
$ws.Range("R9").Value = 0.05490196078431372
$ws.Range("S9").Value = 0.4117647058823529
# Row 10
$ws.Range("B10").Value = 0.1214421252371917
$ws.Range("D10").Value = 0.03542061986084757
$ws.Range("E10").Value = 0.001265022137887413
$ws.Range("F10").Value = 0.058191018342821
$ws.Range("J10").Value = 0.1302972802024036
$ws.Range("O10").Value = 0.02530044275774826
$ws.Range("Q10").Value = 0.228336495888678
$ws.Range("R10").Value = 0.04743833017077799
$ws.Range("S10").Value = 0.3523086654016445
# Row 11
$ws.Range("F11").Value = 0.002710027100271003
$ws.Range("G11").Value = 0.1327913279132791
$ws.Range("J11").Value = 0.08672086720867209
$ws.Range("K11").Value = 0.1517615176151761
$ws.Range("L11").Value = 0.6097560975609756
$ws.Range("S11").Value = 0.01626016260162602
# Row 12
$ws.Range("G12").Value = 0.7287449392712551
$ws.Range("J12").Value = 0.194331983805668
$ws.Range("K12").Value = 0.008097165991902834
$ws.Range("L12").Value = 0.02834008097165992
$ws.Range("S12").Value = 0.04048582995951417
# Row 13
$ws.Range("G13").Value = 0.5555555555555556
$ws.Range("J13").Value = 0.3777777777777778
$ws.Range("S13").Value = 0.06666666666666667
# Row 15
$ws.Range("F15").Value = 0.0308641975308642
$ws.Range("H15").Value = 0.1111111111111111
$ws.Range("I15").Value = 0.07407407407407407
$ws.Range("J15").Value = 0.3858024691358025
$ws.Range("K15").Value = 0.06481481481481481
$ws.Range("M15").Value = 0.009259259259259259
$ws.Range("N15").Value = 0.00308641975308642
$ws.Range("O15").Value = 0.07716049382716049
$ws.Range("S15").Value = 0.2438271604938272
# Row 16
$ws.Range("F16").Value = 0.0421455938697318
$ws.Range("H16").Value = 0.1379310344827586
$ws.Range("I16").Value = 0.103448275862069
$ws.Range("J16").Value = 0.3831417624521073
$ws.Range("K16").Value = 0.1226053639846743
$ws.Range("M16").Value = 0.03831417624521073
$ws.Range("O16").Value = 0.05363984674329502
$ws.Range("S16").Value = 0.1187739463601533
# Row 17
$ws.Range("F17").Value = 0.01923076923076923
$ws.Range("H17").Value = 0.1188811188811189
$ws.Range("I17").Value = 0.1066433566433566
$ws.Range("J17").Value = 0.4493006993006993
$ws.Range("K17").Value = 0.1258741258741259
$ws.Range("M17").Value = 0.008741258741258742
$ws.Range("N17").Value = 0.001748251748251748
$ws.Range("O17").Value = 0.0437062937062937
$ws.Range("S17").Value = 0.1258741258741259
# Row 18
$ws.Range("F18").Value = 0.01418439716312057
$ws.Range("H18").Value = 0.1276595744680851
$ws.Range("I18").Value = 0.09929078014184398
$ws.Range("J18").Value = 0.425531914893617
$ws.Range("K18").Value = 0.09929078014184398
$ws.Range("M18").Value = 0.02836879432624113
$ws.Range("N18").Value = 0.007092198581560284
$ws.Range("O18").Value = 0.09929078014184398
$ws.Range("S18").Value = 0.09929078014184398
# Row 19
$ws.Range("F19").Value = 0.0187793427230047
$ws.Range("H19").Value = 0.1891348088531187
$ws.Range("I19").Value = 0.08517773306505701
$ws.Range("J19").Value = 0.3715627095908786
$ws.Range("K19").Value = 0.1086519114688129
$ws.Range("M19").Value = 0.01676727028839705
$ws.Range("O19").Value = 0.09456740442655935
$ws.Range("S19").Value = 0.1153588195841717
